# Update Rspo2-Znrf3 sheet with new TPM-derived values and add MuSCs-sourced rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 2-6: columns A-D unchanged in content, but some numeric
#     (M..T) values were recomputed with the new TPM data. I..J also changed.
$ws.Range("I2").Value = 0.9891011365778073
$ws.Range("J2").Value = 0.9927075980877177
$ws.Range("M2").Value = 0.7154376666666667
$ws.Range("N2").Value = 2.146313
$ws.Range("O2").Value = 0.05718859144736913
$ws.Range("P2").Value = 0.06498002791311712
$ws.Range("Q2").Value = 1.342286741216778
$ws.Range("R2").Value = 12.080580670951
$ws.Range("S2").Value = 0.05656530079987668
$ws.Range("T2").Value = 0.06450616743330334

$ws.Range("I3").Value = 0.9891011365778073
$ws.Range("J3").Value = 0.9927075980877177
$ws.Range("M3").Value = 1.853215333333333
$ws.Range("N3").Value = 5.559646
$ws.Range("O3").Value = 0.1481369789429594
$ws.Range("P3").Value = 0.1683193235409048
$ws.Range("Q3").Value = 3.476957513493555
$ws.Range("R3").Value = 31.292617621442
$ws.Range("S3").Value = 0.1465224542416838
$ws.Range("T3").Value = 0.167091871384041

$ws.Range("I4").Value = 0.9891011365778073
$ws.Range("J4").Value = 0.9927075980877177
$ws.Range("M4").Value = 2.352629
$ws.Range("N4").Value = 7.057887
$ws.Range("O4").Value = 0.1880576673228452
$ws.Range("P4").Value = 0.2136788503203524
$ws.Range("Q4").Value = 4.413945282494333
$ws.Range("R4").Value = 39.725507542449
$ws.Range("S4").Value = 0.1860080524911973
$ws.Range("T4").Value = 0.2121206182636619

$ws.Range("I5").Value = 0.9891011365778073
$ws.Range("J5").Value = 0.9927075980877177
$ws.Range("M5").Value = 4.500090999999999
$ws.Range("N5").Value = 9.000181999999999
$ws.Range("O5").Value = 0.3597152871109425
$ws.Range("P5").Value = 0.2724821950867064
$ws.Range("Q5").Value = 8.442961231985665
$ws.Range("R5").Value = 50.65776739191399
$ws.Range("S5").Value = 0.3557947993258455
$ws.Range("T5").Value = 0.2704951454061932

$ws.Range("I6").Value = 0.9891011365778073
$ws.Range("J6").Value = 0.9927075980877177
$ws.Range("M6").Value = 3.088773666666667
$ws.Range("N6").Value = 9.266321
$ws.Range("O6").Value = 0.2469014751758839
$ws.Range("P6").Value = 0.2805396031389193
$ws.Range("Q6").Value = 5.795081993240777
$ws.Range("R6").Value = 52.155737939167
$ws.Range("S6").Value = 0.244210529719204
$ws.Range("T6").Value = 0.2784937956005181

# --- New rows 7-11: MuSCs as sending cluster, mirroring target clusters
#     ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Rspo2"
$ws.Range("C7").Value = "Znrf3"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.0206735
$ws.Range("H7").Value = 0.041347
$ws.Range("I7").Value = 0.01089886342219268
$ws.Range("J7").Value = 0.007292401912282354
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7154376666666667
$ws.Range("N7").Value = 2.146313
$ws.Range("O7").Value = 0.05718859144736913
$ws.Range("P7").Value = 0.06498002791311712
$ws.Range("Q7").Value = 0.01479060060183334
$ws.Range("R7").Value = 0.088743603611
$ws.Range("S7").Value = 0.0006232906474924525
$ws.Range("T7").Value = 0.0004738604798137761

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Rspo2"
$ws.Range("C8").Value = "Znrf3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.0206735
$ws.Range("H8").Value = 0.041347
$ws.Range("I8").Value = 0.01089886342219268
$ws.Range("J8").Value = 0.007292401912282354
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.853215333333333
$ws.Range("N8").Value = 5.559646
$ws.Range("O8").Value = 0.1481369789429594
$ws.Range("P8").Value = 0.1683193235409048
$ws.Range("Q8").Value = 0.03831244719366667
$ws.Range("R8").Value = 0.229874683162
$ws.Range("S8").Value = 0.001614524701275547
$ws.Range("T8").Value = 0.001227452156863766

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Rspo2"
$ws.Range("C9").Value = "Znrf3"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.0206735
$ws.Range("H9").Value = 0.041347
$ws.Range("I9").Value = 0.01089886342219268
$ws.Range("J9").Value = 0.007292401912282354
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.352629
$ws.Range("N9").Value = 7.057887
$ws.Range("O9").Value = 0.1880576673228452
$ws.Range("P9").Value = 0.2136788503203524
$ws.Range("Q9").Value = 0.0486370756315
$ws.Range("R9").Value = 0.291822453789
$ws.Range("S9").Value = 0.002049614831647836
$ws.Range("T9").Value = 0.001558232056690432

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Rspo2"
$ws.Range("C10").Value = "Znrf3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.0206735
$ws.Range("H10").Value = 0.041347
$ws.Range("I10").Value = 0.01089886342219268
$ws.Range("J10").Value = 0.007292401912282354
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.500090999999999
$ws.Range("N10").Value = 9.000181999999999
$ws.Range("O10").Value = 0.3597152871109425
$ws.Range("P10").Value = 0.2724821950867064
$ws.Range("Q10").Value = 0.09303263128849999
$ws.Range("R10").Value = 0.372130525154
$ws.Range("S10").Value = 0.003920487785096989
$ws.Range("T10").Value = 0.001987049680513191

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Rspo2"
$ws.Range("C11").Value = "Znrf3"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.0206735
$ws.Range("H11").Value = 0.041347
$ws.Range("I11").Value = 0.01089886342219268
$ws.Range("J11").Value = 0.007292401912282354
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.088773666666667
$ws.Range("N11").Value = 9.266321
$ws.Range("O11").Value = 0.2469014751758839
$ws.Range("P11").Value = 0.2805396031389193
$ws.Range("Q11").Value = 0.06385576239783333
$ws.Range("R11").Value = 0.383134574387
$ws.Range("S11").Value = 0.002690945456679855
$ws.Range("T11").Value = 0.002045807538401188
